$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

function Make-Bold($text) {
    $r = $d.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
    $r.Font.Bold = 1
}

# Title heading text + bold
Replace-Text "デジタル変革の戦略を" "デジタル変革の概要"
Make-Bold "デジタル変革の概要"

Replace-Text "Fabrikam Inc. は、運用効率の向上、カスタマー エクスペリエンスの向上、イノベーションの推進を目的とした包括的なデジタル変革イニシアチブを実施しています。" "Fabrikam Inc. は、運用効率の向上、カスタマー エクスペリエンスの向上、イノベーションの推進を目的とした包括的なデジタル変革イニシアティブに取り組んでいます。"

Replace-Text "次の概要では、現在までに達成された主要な更新プログラムとマイルストーンの概要を示します。" "次の概要では、現在までに達成された主要な更新とマイルストーンについて説明します。"

Replace-Text "キーの更新" "主要な更新"

Replace-Text "オンプレミス アプリケーションの 80% をクラウドに移行する。" "オンプレミス アプリケーションの 80% をクラウドに移行した。"

Replace-Text "スケーラビリティが向上し、IT コストが 25% 削減されました。" "スケーラビリティを向上させ、IT コストを 25% 削減した。"

Replace-Text "業界標準に対するデータ セキュリティとコンプライアンスの強化。" "データ セキュリティを強化し、業界標準に準拠した。"

Replace-Text "意思決定プロセスを効率化するための AI を利用した分析の統合。" "AI を活用した分析を統合し、意思決定プロセスを効率化した。"

Replace-Text "顧客の行動を予測し、マーケティング活動をカスタマイズするための機械学習モデルのデプロイ。" "機械学習モデルをデプロイして顧客の行動を予測し、マーケティング活動をパーソナライズした。"

Replace-Text "手動プロセスの削減により、生産性が 30% 向上します。" "手動プロセスを削減した結果、生産性が 30% 向上した。"

Replace-Text "セルフサービス機能を備えた新しい顧客ポータルを起動します。" "セルフサービス機能を備えた新しい顧客ポータルを立ち上げた。"

Replace-Text "24 時間 365 日のカスタマー サポートのためのチャットボットの導入により、応答時間を 50% 短縮します。" "チャットボットを導入して 24 時間 365 日のカスタマー サポートを実現し、応答時間を 50% 短縮した。"

Replace-Text "過去 1 年間で顧客満足度スコアが 20% 向上しました。" "過去 1 年間で顧客満足度スコアが 20% 向上した。"

# "プロセスの自動化" heading: bold only, no text change
Make-Bold "プロセスの自動化"

Replace-Text "ルーチン タスク用のロボット プロセス自動化 (RPA) の実装。" "ロボティック プロセス オートメーション (RPA) を実装し、ルーチン タスクを自動化した。"

Replace-Text "主要業務の処理時間を40%削減。" "主要業務の処理時間を 40% 削減した。"

Replace-Text "組織内のより戦略的な役割への人事の再割り当て。" "組織内でより戦略的な役割に人的資源を再配置した。"

Replace-Text "全従業員を対象にデジタルリテラシープログラムを実施。" "全従業員を対象にデジタル リテラシー プログラムを実施した。"

Replace-Text "新しいテクノロジに関するコースを含む新しい e ラーニング プラットフォームを立ち上げました。" "新しいテクノロジに関するコースを含む新しい e ラーニング プラットフォームを立ち上げた。"

Replace-Text "従業員のエンゲージメントと新しいツールの導入を 35% 増やしました。" "従業員のエンゲージメントと新しいツールの導入率を 35% 増加させた。"

Replace-Text "Q1 2024: クラウド インフラストラクチャへの移行が完了しました。" "2024 年第 1 四半期: クラウド インフラストラクチャへの移行が完了した。"

Replace-Text "2024 年第 2 四半期: AI 主導の分析プラットフォームを開始しました。" "2024 年第 2 四半期: AI 駆動型の分析プラットフォームを立ち上げた。"

Replace-Text "Q3 2024: 新しいデジタル カスタマー ポータルが導入されました。" "2024 年第 3 四半期: 新しいデジタル カスタマー ポータルを導入した。"

Replace-Text "Q4 2024: ルーチンプロセスの50%自動化を達成しました。" "2024 年第 4 四半期: ルーチン プロセスの自動化率が 50% に達した。"

# "今後の予定" heading: bold only, no text change
Make-Bold "今後の予定"

Replace-Text "新しい機能とサービスを使用して、デジタル カスタマー エクスペリエンスをさらに強化します。" "新しい機能とサービスを活用して、デジタル カスタマー エクスペリエンスをさらに強化します。"

Replace-Text "進化する脅威から保護するためのサイバーセキュリティ対策に重点を置く。" "サイバーセキュリティ対策に重点を置き、進化する脅威からの保護を強化します。"

Replace-Text "Fabrikam Inc.のデジタル変革の取り組みにより、効率性、顧客満足度、全体的なビジネス パフォーマンスが大幅に向上しました。" "Fabrikam Inc. はデジタル変革に取り組むことで、効率性、顧客満足度、全体的なビジネス パフォーマンスを大幅に改善しました。"

Replace-Text "組織は、将来の成長とイノベーションを促進するためにテクノロジを活用することに引き続き取り組んでいます。" "同社は、今後もテクノロジーを活用して将来的な成長とイノベーションの促進に取り組み続けていく方針です。"

Write-Output "All replacements applied."
